# Pipeline edit: remove the obsolete "Test Bundle 3 Provider" deal row from
# the Deals_data sheet (it now only ships as a row on the Sheet1 tab),
# clearing the stale Courier-New formatting that had been applied to the
# remaining deal-name/provider cells, and refresh the saved selections.

$wb = $excel.ActiveWorkbook

# --- Deals_data: drop row 4 (VAS / Test Bundle 3 Provider / ... / Cape Town) ---
$deals = $wb.Worksheets.Item("Deals_data")

# The provider/deal-name cells in the rows that remain (B2:C3) had picked up
# a Courier New style; clear it back to the sheet default.
$deals.Range("B2:C3").ClearFormats()

# Remove the whole 4th row, shifting nothing else (it was the last row).
$deals.Rows.Item(4).Delete()

# Leave the cursor parked where the pipeline run left it.
[void]$deals.Activate()
[void]$deals.Range("C8").Select()

# --- Sheet1 (3rd tab): just a saved-selection refresh, content unchanged ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
[void]$sheet1.Activate()
[void]$sheet1.Range("A3:E5").Select()

# Restore focus to the Deals_data tab, which is the one marked active/selected.
[void]$deals.Activate()
